
# Werkelijk resultaat kolom toegevoegd aan overview.
# Inserts a new "Werkelijk resultaat" column right after "Verwacht resultaat"
# (column L), shifting every following column one position to the right,
# and keeps the AutoFilter / filter-database defined name in sync with the
# now one-column-wider header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overzicht")

# Remember the width of the donor column (L) before the insert so the new
# column (M) can be sized the same way.
$donorWidth = $ws.Columns("L:L").ColumnWidth

# Insert a new column before the current column M ("Niet toegewezen regels").
# Shift-right insert copies the formatting of the column immediately to the
# left (L), so the new header cell automatically picks up the bold/blue
# header style already used across the row.
$ws.Columns("M:M").Insert()

# New header text -> lands in the freshly inserted M1 cell.
$ws.Range("M1").Value = "Werkelijk resultaat"

# Match the column width of the donor column as closely as possible.
$ws.Columns("M:M").ColumnWidth = $donorWidth

# The autofilter used to span A1:U1; after inserting the extra column it
# must cover one more column, A1:V1. Toggle it off and back on over the new
# range (re-enabling over the same range is a no-op, so the off/on pair is
# required to move the boundary).
$ws.AutoFilterMode = $false
$ws.Range("A1:V1").AutoFilter()

# Keep the hidden _FilterDatabase defined name lined up with the new
# autofilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Overzicht!_FilterDatabase") {
        $n.RefersTo = "=Overzicht!`$A`$1:`$V`$1"
    }
}

# Restore the selection to the cell the author ended up on after adding the
# column.
$ws.Range("M8").Select()
